$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '26.698.18'
$ws.Range("E2").Value = '  +0.35%  '

Set-TextValue "D3" '1.599.69'
$ws.Range("E3").Value = '  +0.25%  '

$ws.Range("E4").Value = '  +0.20%  '

Set-TextValue "D5" '211.48'
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("E6").Value = '  -0.38%  '

$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("E8").Value = '  +0.49%  '

$ws.Range("E9").Value = '  +1.28%  '

Set-TextValue "D10" '19.55'
$ws.Range("E10").Value = '  +0.49%  '

Set-TextValue "D11" '0.0841'
$ws.Range("E11").Value = '  +0.77%  '

Set-TextValue "D12" '1.824.33'

Set-TextValue "D13" '1.589.47'
$ws.Range("E13").Value = '  -1.61%  '

Set-TextValue "D14" '4.04'
$ws.Range("E14").Value = '  +0.63%  '

Set-TextValue "D16" '65.32'
$ws.Range("E16").Value = '  +1.42%  '

Set-TextValue "D17" '26.681.75'
$ws.Range("E17").Value = '  +0.32%  '

Set-TextValue "D18" '0.0₃0757'
$ws.Range("E18").Value = '  +3.60%  '

Set-TextValue "D19" '209.89'
$ws.Range("E19").Value = '  +0.93%  '

$ws.Range("E20").Value = '  +0.18%  '

$ws.Range("E21").Value = '  +4.06%  '

Set-TextValue "D22" '4.28'
$ws.Range("E22").Value = '  +0.77%  '

Set-TextValue "D23" '2.31'
$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("E24").Value = '  +0.88%  '

Set-TextValue "D25" '142.87'
$ws.Range("E25").Value = '  -1.67%  '

Set-TextValue "D26" '1.00'
$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("E28").Value = '  +0.38%  '

$ws.Range("E29").Value = '  +0.81%  '

$ws.Range("E30").Value = '  +2.66%  '

$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("E32").Value = '  +0.95%  '

$ws.Range("E33").Value = '  +1.71%  '

Set-TextValue "D34" '1.290.84'
$ws.Range("E34").Value = '  +0.71%  '

Set-TextValue "D35" '0.618'
$ws.Range("E35").Value = '  -5.13%  '

Set-TextValue "D36" '2.47'
$ws.Range("E36").Value = '  +0.83%  '

$ws.Range("E37").Value = '  +0.39%  '

$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("E39").Value = '  +15.86%  '

Set-TextValue "D40" '0.824'
$ws.Range("E40").Value = '  -2.29%  '

$ws.Range("E41").Value = '  -0.73%  '

Set-TextValue "D42" '0.786'
$ws.Range("E42").Value = '  +0.23%  '

$ws.Range("E43").Value = '  -0.72%  '

Set-TextValue "D44" '63.15'
$ws.Range("E44").Value = '  -0.99%  '

Set-TextValue "D45" '1.735.76'
$ws.Range("E45").Value = '  +0.21%  '

Set-TextValue "D46" '91.08'
$ws.Range("E46").Value = '  +1.74%  '

$ws.Range("E47").Value = '  -0.58%  '

$ws.Range("E48").Value = '  -1.58%  '

$ws.Range("E49").Value = '  +0.56%  '

$ws.Range("E50").Value = '  +0.13%  '

Set-TextValue "D51" '7.35'
$ws.Range("E51").Value = '  -1.06%  '
